$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new most-recent-quarter + prior quarter),
# shifting the existing quarterly data (old D:K) right to F:M.
$ws.Range("D:E").Insert()

# The newly inserted D:E columns pick up the plain default column style; restore the
# correct number formats / fonts / alignment by copying formats from column F (which
# now holds what used to be column D) across the full data block.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 37 and 79 are section-header rows (only column B is used) - the blanket
# format paste above should not leave stray empty cells in D:E there.
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()

# Match the new columns' width to their neighboring (bestFit) quarterly column.
$ws.Columns("D:E").ColumnWidth = $ws.Columns("F").ColumnWidth

# Populate the two new columns with the new quarter figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 56300
$ws.Range("E8").Value = 35300
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = -2000
$ws.Range("E15").Value = -700
$ws.Range("D17").Value = 13400
$ws.Range("E17").Value = 7300
$ws.Range("D18").Value = 42900
$ws.Range("E18").Value = 28000
$ws.Range("D20").Value = -26700
$ws.Range("E20").Value = -17200
$ws.Range("D21").Value = 18200
$ws.Range("E21").Value = 11500
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 16200
$ws.Range("E23").Value = 10800
$ws.Range("D24").Value = 3000
$ws.Range("E24").Value = 1900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 13200
$ws.Range("E26").Value = 8900
$ws.Range("D27").Value = 13200
$ws.Range("E27").Value = 8900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 26700
$ws.Range("E32").Value = 17200
$ws.Range("D33").Value = 13200
$ws.Range("E33").Value = 8900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 13200
$ws.Range("E35").Value = 8900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 118800
$ws.Range("E41").Value = 145200
$ws.Range("D42").Value = 161100
$ws.Range("E42").Value = 58100
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 41700
$ws.Range("E48").Value = 19000
$ws.Range("D49").Value = 249700
$ws.Range("E49").Value = 42100
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 4655200
$ws.Range("E54").Value = 3035500
$ws.Range("D57").Value = 2800
$ws.Range("E57").Value = 2100
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 274400
$ws.Range("E61").Value = 260400
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 3952300
$ws.Range("E66").Value = 2707400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 112100
$ws.Range("E72").Value = 99000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 703000
$ws.Range("E76").Value = 328100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 13200
$ws.Range("E81").Value = 8900
$ws.Range("D83").Value = 2000
$ws.Range("E83").Value = 700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 13600
$ws.Range("E89").Value = 14200
$ws.Range("D91").Value = -1600
$ws.Range("E91").Value = -400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 143100
$ws.Range("E94").Value = -80200
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -79200
$ws.Range("E100").Value = 56800
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 77500
$ws.Range("E102").Value = -9200
